# Revert "removing _new.xlsx in the databases"
#
# The SUPPLY.xlsx workbook had previously gained three extra columns on the
# HEATING and COOLING sheets ("primary_components", "secondary_components",
# "tertiary_components") that listed component codes such as CH1/CH2/CT1/CT2/
# AC1/HEX1 (or "-" where not applicable). This commit reverts that change,
# deleting those three columns again from both sheets so the tables go back
# to: Description, code, (system,) feedstock, scale, efficiency,
# CAPEX_USD2015kW, LT_yr, O&M_%, IR_%, reference
# - matching the HOT_WATER / ELECTRICITY sheets, which were never touched.

$wb = $excel.ActiveWorkbook

# HEATING sheet: columns C:E are primary_components/secondary_components/
# tertiary_components. Remove them so C becomes feedstock again.
$wsHeating = $wb.Worksheets.Item("HEATING")
$wsHeating.Range("C1:E1").EntireColumn.Delete()

# COOLING sheet: same three extra columns (C:E), here followed by "system".
$wsCooling = $wb.Worksheets.Item("COOLING")
$wsCooling.Range("C1:E1").EntireColumn.Delete()
